$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "309.76"
Set-TextValue "E2" "-2.48%"
Set-TextValue "G2" "17"
Set-TextValue "D3" "37.70"
Set-TextValue "E3" "-4.95%"
Set-TextValue "G3" "17"
Set-TextValue "D4" "5.101"
Set-TextValue "E4" "-0.78%"
Set-TextValue "G4" "17"
Set-TextValue "D5" "0.07857"
Set-TextValue "E5" "-4.33%"
Set-TextValue "G5" "17"
Set-TextValue "D6" "1.969"
Set-TextValue "E6" "-3.09%"
Set-TextValue "G6" "17"
Set-TextValue "D7" "4.372"
Set-TextValue "E7" "1.99%"
Set-TextValue "G7" "17"
Set-TextValue "D8" "8.292"
Set-TextValue "E8" "0.04%"
Set-TextValue "G8" "17"
Set-TextValue "D9" "3.133"
Set-TextValue "E9" "-1.95%"
Set-TextValue "G9" "17"
Set-TextValue "D10" "0.9276"
Set-TextValue "E10" "-0.45%"
Set-TextValue "G10" "17"
Set-TextValue "D11" "0.1353"
Set-TextValue "E11" "-4.31%"
Set-TextValue "G11" "17"
Set-TextValue "D12" "0.2006"
Set-TextValue "E12" "0.42%"
Set-TextValue "G12" "17"
Set-TextValue "D13" "0.08937"
Set-TextValue "E13" "-1.09%"
Set-TextValue "G13" "17"
Set-TextValue "D14" "0.03468"
Set-TextValue "E14" "-0.27%"
Set-TextValue "G14" "17"
Set-TextValue "D15" "0.09741"
Set-TextValue "E15" "-0.61%"
Set-TextValue "G15" "17"
Set-TextValue "D16" "0.001393"
Set-TextValue "E16" "-0.80%"
Set-TextValue "G16" "17"
Set-TextValue "D17" "0.005939"
Set-TextValue "E17" "-3.43%"
Set-TextValue "G17" "17"
Set-TextValue "E18" "1,777.67%"
Set-TextValue "G18" "17"
Set-TextValue "E19" "-2.52%"
Set-TextValue "G19" "17"
Set-TextValue "D20" "0.3466"
Set-TextValue "E20" "-0.16%"
Set-TextValue "G20" "17"
Set-TextValue "D21" "0.1294"
Set-TextValue "E21" "0.26%"
Set-TextValue "G21" "17"
Set-TextValue "D22" "5.010"
Set-TextValue "E22" "2.19%"
Set-TextValue "G22" "17"
Set-TextValue "D23" "0.2514"
Set-TextValue "E23" "2.60%"
Set-TextValue "G23" "17"
Set-TextValue "D24" "0.04321"
Set-TextValue "E24" "-0.07%"
Set-TextValue "G24" "17"
Set-TextValue "D25" "0.001219"
Set-TextValue "E25" "-0.49%"
Set-TextValue "G25" "17"
Set-TextValue "D26" "0.004540"
Set-TextValue "E26" "-4.80%"
Set-TextValue "G26" "17"
Set-TextValue "D27" "0.0001352"
Set-TextValue "E27" "3.98%"
Set-TextValue "G27" "17"
Set-TextValue "G28" "17"
Set-TextValue "G29" "17"
Set-TextValue "G30" "17"
Set-TextValue "G31" "17"
Set-TextValue "G32" "17"
Set-TextValue "G33" "17"
Set-TextValue "G34" "17"
Set-TextValue "G35" "17"
Set-TextValue "G36" "17"
Set-TextValue "G37" "17"
Set-TextValue "G38" "17"
Set-TextValue "D39" "0.02301"
Set-TextValue "E39" "3.73%"
Set-TextValue "G39" "17"
Set-TextValue "D40" "0.05057"
Set-TextValue "E40" "-3.15%"
Set-TextValue "G40" "17"
Set-TextValue "D41" "0.007470"
Set-TextValue "E41" "-0.53%"
Set-TextValue "G41" "17"
Set-TextValue "D42" "0.009861"
Set-TextValue "E42" "0.41%"
Set-TextValue "G42" "17"
Set-TextValue "D43" "0.1356"
Set-TextValue "E43" "-1.69%"
Set-TextValue "G43" "17"
Set-TextValue "D44" "0.001983"
Set-TextValue "E44" "-7.81%"
Set-TextValue "G44" "17"
Set-TextValue "D45" "0.008772"
Set-TextValue "E45" "-10.90%"
Set-TextValue "G45" "17"
Set-TextValue "D46" "0.00006803"
Set-TextValue "E46" "3.13%"
Set-TextValue "G46" "17"
Set-TextValue "E47" "0.22%"
Set-TextValue "G47" "17"
Set-TextValue "D48" "0.003004"
Set-TextValue "E48" "8.72%"
Set-TextValue "G48" "17"
Set-TextValue "E49" "8.55%"
Set-TextValue "G49" "17"
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "0.22%"
Set-TextValue "G50" "17"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "0.22%"
Set-TextValue "G51" "17"
